# Apply the "updated 4.0 files and mdl" edit:
#  - About sheet: bump the last-modified date in C1
#  - MCF sheet: raise several plant-type capacity factors to 1 (100%)
#  - MCF sheet: update the saved cell selection to B17

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$mcf   = $wb.Worksheets.Item("MCF")

# --- About sheet: last-modified date (C1) 1/29/2024 -> 4/10/2024 ---
$about.Range("C1").Value = "4/10/2024"

# --- MCF sheet: capacity factors bumped to 1 ---
$mcf.Range("B2").Value  = 1
$mcf.Range("B3").Value  = 1
$mcf.Range("B4").Value  = 1
$mcf.Range("B6").Value  = 1
$mcf.Range("B10").Value = 1
$mcf.Range("B11").Value = 1
$mcf.Range("B12").Value = 1
$mcf.Range("B13").Value = 1
$mcf.Range("B14").Value = 1
$mcf.Range("B16").Value = 1
$mcf.Range("B17").Value = 1
$mcf.Range("B18").Value = 1

# --- MCF sheet: update the active selection shown when the sheet is opened ---
$mcf.Activate()
$mcf.Range("B17").Select()
